$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 2: remove the duplicated "Play All Star Knockout Free..."
# bold paragraph that used to sit right before the closing italic
# paragraph near the end of the document.
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text
    if ($i -ne 1 -and $paraText -eq "Play All Star Knockout Free - Exciting medium-low variance fruit slot`r") {
        $para.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# Change 3: swap the closing italic paragraph's text for the new
# image-generation prompt. Done before Change 1 below so that this
# Find only ever matches the single occurrence at the end of the
# document (Change 1 introduces a second, unrelated copy of this
# phrase inside the new Meta description paragraph).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Try All Star Knockout for free with its Random Bonus feature and 5 Multipliers. With a bet range of " + [char]0x20AC + "0.20 to " + [char]0x20AC + "400.00, it's suitable for all players.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Prompt: Create a cartoon-style feature image for All Star Knockout that features a happy Maya warrior with glasses. The image should be colorful and eye-catching, with the Maya warrior front and center. The warrior should be wearing a headpiece with feathers and a pair of glasses, imparting a modern twist to their traditional attire. The background should include elements that reference the game, such as stars, fruits, and the prize board. The overall feel of the image should be playful and inviting, encouraging players to try out the game.",
    2
)

# ------------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the
# very first paragraph (the Heading1 title).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$insertionPoint = $metaPara.Range.Start
$metaRange = $d.Range($insertionPoint, $insertionPoint)
$metaRange.InsertAfter("Meta description: Try All Star Knockout for free with its Random Bonus feature and 5 Multipliers. With a bet range of " + [char]0x20AC + "0.20 to " + [char]0x20AC + "400.00, it's suitable for all players.")

# Bold just the "Meta description" label (16 characters).
$labelRange = $d.Range($insertionPoint, $insertionPoint + 16)
$labelRange.Font.Bold = 1
